$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("C1").Value = "PPT height"
$ws.Range("D1").Value = "PPT width"

# New data values (PPT height/width in inches, based on 96 DPI)
$ws.Range("C2").Value = 7.03
$ws.Range("D2").Value = 12.5

$ws.Range("C3").Value = 9.375
$ws.Range("D3").Value = 12.5

$ws.Range("C4").Value = 12.5
$ws.Range("D4").Value = 12.5

# Update selection to match target state
$ws.Range("A6:XFD7").Select()
